$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "Meta description" paragraph that sat right under
#    the H1 title ("Play Fruit Slot Free: Unique Board-Style Layout").
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold-only paragraph reading
#    "Play Fruit Slot Free: Unique Board-Style Layout" right before the
#    final ("Create a feature image...") paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)

$xmlFrag = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Slot Free: Unique Board-Style Layout</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xmlFrag)

# ------------------------------------------------------------------
# 3) Replace the old image-prompt paragraph's text with the new meta
#    description copy, keeping its existing (italic) run formatting.
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "Create a feature image that captures the fun and excitement of Fruit Slot! Your image should be in a cartoon style and should feature a happy Maya warrior wearing glasses, surrounded by colorful fruit symbols. Be creative and playful with your design, incorporating the game's Asian arcade theme. Consider including the ring pattern of symbols in your design, as well as some of the potential multipliers that players can win. Your image should be eye-catching and convey the game's unique twist on traditional slots."
$find.Replacement.Text = "Experience a unique gaming experience on Fruit Slot, with a board-style layout and symbol betting. Play now for free."
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
